$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.543.11"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.815.73"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +24.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.296"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0675"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("D12").Value = "2.081.61"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "1.823.04"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.639"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "34.532.06"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "0.0₃0779"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.42%  "
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0520"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "88.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.658"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "15.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.46%  "
$ws.Range("D38").Value = "1.320.08"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0190"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.956"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  +4.33%  "
$ws.Range("D47").Value = "1.982.86"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0611"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
